$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.034.15'
$ws.Cells.Item(2, 5).Value = '  -3.29%  '

$ws.Cells.Item(3, 4).Value = '3.161.92'
$ws.Cells.Item(3, 5).Value = '  -8.41%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).Value = '562.70'
$ws.Cells.Item(5, 5).Value = '  -4.11%  '

$ws.Cells.Item(6, 4).Value = '169.80'
$ws.Cells.Item(6, 5).Value = '  -3.64%  '

$ws.Cells.Item(7, 2).Value = 'XRP'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Cells.Item(7, 4).Value = '0.610'
$ws.Cells.Item(7, 5).Value = '  +1.07%  '

$ws.Cells.Item(8, 2).Value = 'USDC'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(8, 4).Value = '1.00'
$ws.Cells.Item(8, 5).Value = '  +0.11%  '

$ws.Cells.Item(9, 4).Value = '3.155.20'
$ws.Cells.Item(9, 5).Value = '  -8.56%  '

$ws.Cells.Item(10, 5).Value = '  -6.50%  '

$ws.Cells.Item(11, 5).Value = '  -5.06%  '

$ws.Cells.Item(12, 4).Value = '0.396'
$ws.Cells.Item(12, 5).Value = '  -5.30%  '

$ws.Cells.Item(13, 4).Value = '3.715.93'
$ws.Cells.Item(13, 5).Value = '  -8.30%  '

$ws.Cells.Item(14, 5).Value = '  +1.11%  '

$ws.Cells.Item(15, 4).Value = '27.26'
$ws.Cells.Item(15, 5).Value = '  -7.42%  '

$ws.Cells.Item(16, 4).Value = '64.034.19'
$ws.Cells.Item(16, 5).Value = '  -3.16%  '

$ws.Cells.Item(17, 4).Value = '0.0000163'
$ws.Cells.Item(17, 5).Value = '  -5.64%  '

$ws.Cells.Item(18, 4).Value = '3.170.68'
$ws.Cells.Item(18, 5).Value = '  -8.24%  '

$ws.Cells.Item(19, 4).Value = '5.72'
$ws.Cells.Item(19, 5).Value = '  -3.84%  '

$ws.Cells.Item(20, 4).Value = '13.01'
$ws.Cells.Item(20, 5).Value = '  -5.54%  '

$ws.Cells.Item(21, 4).Value = '354.11'
$ws.Cells.Item(21, 5).Value = '  -5.03%  '

$ws.Cells.Item(22, 4).Value = '7.20'
$ws.Cells.Item(22, 5).Value = '  -5.55%  '

$ws.Cells.Item(23, 4).Value = '0.998'
$ws.Cells.Item(23, 5).Value = '  +0.07%  '

$ws.Cells.Item(24, 4).Value = '69.12'
$ws.Cells.Item(24, 5).Value = '  -5.49%  '

$ws.Cells.Item(25, 4).Value = '0.503'
$ws.Cells.Item(25, 5).Value = '  -6.49%  '

$ws.Cells.Item(26, 4).Value = '0.0000118'
$ws.Cells.Item(26, 5).Value = '  -4.90%  '

$ws.Cells.Item(27, 4).Value = '9.63'
$ws.Cells.Item(27, 5).Value = '  -1.49%  '

$ws.Cells.Item(28, 5).Value = '  -1.62%  '

$ws.Cells.Item(29, 5).Value = '  -0.03%  '

$ws.Cells.Item(30, 4).Value = '5.62'
$ws.Cells.Item(30, 5).Value = '  -3.73%  '

$ws.Cells.Item(31, 4).Value = '0.999'
$ws.Cells.Item(31, 5).Value = '  -0.09%  '

$ws.Cells.Item(32, 5).Value = '  -5.23%  '

$ws.Cells.Item(33, 5).Value = '  -6.70%  '

$ws.Cells.Item(34, 4).Value = '6.63'
$ws.Cells.Item(34, 5).Value = '  -5.98%  '

$ws.Cells.Item(35, 5).Value = '  -5.65%  '

$ws.Cells.Item(36, 4).Value = '1.43'
$ws.Cells.Item(36, 5).Value = '  -8.01%  '

$ws.Cells.Item(37, 4).Value = '154.92'
$ws.Cells.Item(37, 5).Value = '  -4.19%  '

$ws.Cells.Item(38, 4).Value = '0.807'
$ws.Cells.Item(38, 5).Value = '  -8.63%  '

$ws.Cells.Item(39, 4).Value = '25.83'
$ws.Cells.Item(39, 5).Value = '  -8.99%  '

$ws.Cells.Item(40, 4).Value = '2.55'
$ws.Cells.Item(40, 5).Value = '  -3.09%  '

$ws.Cells.Item(41, 5).Value = '  -5.97%  '

$ws.Cells.Item(42, 4).Value = '2.615.54'
$ws.Cells.Item(42, 5).Value = '  -5.79%  '

$ws.Cells.Item(43, 5).Value = '  -7.09%  '

$ws.Cells.Item(44, 4).Value = '6.00'
$ws.Cells.Item(44, 5).Value = '  -7.12%  '

$ws.Cells.Item(45, 2).Value = 'OKB'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(45, 4).Value = '39.46'
$ws.Cells.Item(45, 5).Value = '  -1.11%  '

$ws.Cells.Item(46, 2).Value = 'Hedera'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(46, 4).Value = '0.0658'
$ws.Cells.Item(46, 5).Value = '  -4.65%  '

$ws.Cells.Item(47, 4).Value = '326.74'
$ws.Cells.Item(47, 5).Value = '  -2.78%  '

$ws.Cells.Item(48, 4).Value = '23.86'
$ws.Cells.Item(48, 5).Value = '  -5.38%  '

$ws.Cells.Item(49, 4).Value = '0.0270'
$ws.Cells.Item(49, 5).Value = '  -7.57%  '

$ws.Cells.Item(50, 5).Value = '  -0.54%  '

$ws.Cells.Item(51, 4).Value = '0.999'
$ws.Cells.Item(51, 5).Value = '  -0.07%  '
